# Agricultural Products Net Import Value (Import-Export): F3=(F3.1-F3.2)
#
# Source sheet (table of feature definitions): the per-capita / F3.3
# normalisation is dropped from feature F3 - both its long name and its
# formula description shed the "per Capita" / "/F3.3" suffix. The F1.A
# label is renamed F1*. Data sheet column G (F3) is recomputed from the
# raw F3.1/F3.2 columns (D-E) instead of the old per-capita ratio.

$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item("Data")
$wsSource = $wb.Worksheets.Item("Source")

# ---------------------------------------------------------------------
# 1) Source sheet text edits
# ---------------------------------------------------------------------
$wsSource.Range("A3").Value = "F1*"
$wsSource.Range("D8").Value = "F3=(F3.1-F3.2)"
$wsSource.Range("B8").Value = "Agricultural Products Net Import Value (Import-Export)"

# ---------------------------------------------------------------------
# 2) Data sheet: column G (F3) becomes D-E instead of a stored ratio
# ---------------------------------------------------------------------
$wsData.Range("G2:G37").Formula = "=D2-E2"
$wsData.Range("G2:G37").NumberFormat = "0"

# ---------------------------------------------------------------------
# 3) Cosmetic re-highlight of column F (F3.3) to match the reworked
#    F3 derivation, and the matching Source-sheet row (F3.3 / Population).
# ---------------------------------------------------------------------
$accentRange = $wsData.Range("F2:F37")
$accentRange.Font.ThemeColor = 6
$accentRange.Font.TintAndShade = -0.249977111117893

$wsData.Range("F1").Font.Bold = $true
$wsData.Range("F1").Font.ThemeColor = 6
$wsData.Range("F1").Font.TintAndShade = -0.249977111117893

$sourceRow7 = $wsSource.Range("A7:C7")
$sourceRow7.Font.ThemeColor = 6
$sourceRow7.Font.TintAndShade = -0.249977111117893

$wsSource.Range("D7").Font.Underline = $true
$wsSource.Range("D7").Font.ThemeColor = 6
$wsSource.Range("D7").Font.TintAndShade = -0.249977111117893

# ---------------------------------------------------------------------
# 4) View state: Source becomes the active/selected tab
# ---------------------------------------------------------------------
$wsSource.Activate()
$wsData.Range("F4").Select()
